$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    14 = '2 ч. 06 мин. 12 сек.'
    16 = '16 ч. 04 мин. 40 сек.'
    22 = '5 ч. 07 мин. 34 сек.'
    24 = '1 ч. 30 мин. 02 сек.'
    25 = '1 ч. 24 мин. 08 сек.'
    26 = '2 ч. 01 мин. 54 сек.'
    29 = '1 ч. 18 мин. 03 сек.'
    30 = '9 ч. 16 мин. 09 сек.'
    32 = '1 ч. 37 мин. 02 сек.'
    36 = '1 ч. 03 мин. 45 сек.'
    42 = '3 ч. 45 мин. 01 сек.'
    43 = '3 ч. 05 мин. 14 сек.'
    44 = '0 ч. 57 мин. 09 сек.'
    47 = '3 ч. 08 мин. 15 сек.'
    50 = '5 ч. 09 мин. 12 сек.'
    52 = '3 ч. 20 мин. 07 сек.'
    59 = '3 ч. 40 мин. 05 сек.'
    62 = '4 ч. 34 мин. 09 сек.'
    66 = '4 ч. 05 мин. 34 сек.'
    67 = '10 ч. 28 мин. 03 сек.'
    73 = '4 ч. 04 мин. 45 сек.'
    79 = '6 ч. 59 мин. 05 сек.'
    92 = '2 ч. 02 мин. 20 сек.'
    94 = '2 ч. 01 мин. 36 сек.'
    97 = '1 ч. 21 мин. 02 сек.'
    110 = '9 ч. 41 мин. 05 сек.'
    111 = '15 ч. 48 мин. 01 сек.'
    113 = '9 ч. 02 мин. 57 сек.'
    117 = '12 ч. 29 мин. 05 сек.'
    118 = '2 ч. 30 мин. 00 сек.'
    121 = '2 ч. 00 мин. 00 сек.'
    125 = '6 ч. 15 мин. 07 сек.'
    127 = '1 ч. 29 мин. 00 сек.'
    128 = '3 ч. 08 мин. 29 сек.'
    140 = '12 ч. 08 мин. 12 сек.'
    146 = '6 ч. 58 мин. 01 сек.'
    147 = '1 ч. 27 мин. 06 сек.'
    148 = '1 ч. 53 мин. 01 сек.'
    149 = '2 ч. 45 мин. 03 сек.'
    155 = '6 ч. 05 мин. 50 сек.'
    160 = '8 ч. 00 мин. 12 сек.'
    163 = '3 ч. 02 мин. 05 сек.'
    165 = '7 ч. 32 мин. 00 сек.'
    169 = '6 ч. 35 мин. 00 сек.'
    172 = '1 ч. 25 мин. 03 сек.'
    175 = '2 ч. 24 мин. 02 сек.'
    179 = '3 ч. 54 мин. 00 сек.'
    188 = '1 ч. 59 мин. 04 сек.'
    192 = '12 ч. 34 мин. 01 сек.'
    193 = '1 ч. 16 мин. 09 сек.'
    199 = '2 ч. 05 мин. 49 сек.'
    211 = '14 ч. 44 мин. 05 сек.'
    213 = '3 ч. 09 мин. 52 сек.'
    218 = '3 ч. 09 мин. 45 сек.'
    223 = '1 ч. 57 мин. 06 сек.'
    229 = '5 ч. 00 мин. 41 сек.'
    232 = '2 ч. 07 мин. 30 сек.'
    237 = '2 ч. 07 мин. 58 сек.'
    239 = '3 ч. 10 мин. 07 сек.'
    240 = '5 ч. 28 мин. 06 сек.'
    242 = '5 ч. 01 мин. 28 сек.'
    243 = '9 ч. 01 мин. 24 сек.'
    259 = '2 ч. 10 мин. 09 сек.'
    267 = '1 ч. 55 мин. 03 сек.'
    268 = '3 ч. 07 мин. 24 сек.'
    270 = '2 ч. 01 мин. 45 сек.'
    276 = '3 ч. 07 мин. 06 сек.'
    277 = '4 ч. 16 мин. 04 сек.'
    279 = '10 ч. 43 мин. 06 сек.'
    282 = '3 ч. 20 мин. 00 сек.'
    297 = '9 ч. 08 мин. 30 сек.'
    299 = '2 ч. 07 мин. 00 сек.'
    302 = '6 ч. 09 мин. 36 сек.'
    303 = '8 ч. 08 мин. 36 сек.'
    304 = '4 ч. 37 мин. 01 сек.'
    306 = '1 ч. 59 мин. 09 сек.'
    308 = '2 ч. 01 мин. 12 сек.'
    310 = '3 ч. 05 мин. 08 сек.'
    316 = '16 ч. 03 мин. 07 сек.'
    318 = '2 ч. 10 мин. 02 сек.'
    320 = '12 ч. 09 мин. 40 сек.'
    324 = '2 ч. 02 мин. 38 сек.'
    329 = '1 ч. 57 мин. 06 сек.'
    330 = '4 ч. 22 мин. 08 сек.'
    335 = '5 ч. 07 мин. 29 сек.'
    338 = '3 ч. 57 мин. 06 сек.'
    341 = '2 ч. 08 мин. 37 сек.'
    343 = '2 ч. 12 мин. 05 сек.'
    350 = '2 ч. 23 мин. 09 сек.'
    351 = '4 ч. 29 мин. 08 сек.'
    352 = '5 ч. 08 мин. 18 сек.'
    354 = '1 ч. 40 мин. 08 сек.'
    355 = '3 ч. 19 мин. 01 сек.'
    356 = '9 ч. 23 мин. 02 сек.'
    357 = '2 ч. 12 мин. 08 сек.'
    364 = '12 ч. 08 мин. 25 сек.'
    365 = '10 ч. 10 мин. 08 сек.'
    370 = '2 ч. 05 мин. 13 сек.'
    373 = '2 ч. 02 мин. 48 сек.'
    378 = '10 ч. 00 мин. 05 сек.'
    380 = '3 ч. 08 мин. 03 сек.'
    385 = '1 ч. 42 мин. 04 сек.'
    386 = '2 ч. 07 мин. 04 сек.'
    387 = '10 ч. 00 мин. 25 сек.'
    395 = '2 ч. 00 мин. 07 сек.'
    397 = '3 ч. 12 мин. 01 сек.'
    398 = '2 ч. 04 мин. 25 сек.'
    400 = '2 ч. 10 мин. 00 сек.'
    406 = '4 ч. 09 мин. 16 сек.'
    407 = '5 ч. 13 мин. 04 сек.'
    408 = '2 ч. 46 мин. 09 сек.'
    410 = '7 ч. 17 мин. 03 сек.'
    411 = '1 ч. 57 мин. 06 сек.'
    412 = '2 ч. 09 мин. 38 сек.'
    414 = '2 ч. 02 мин. 21 сек.'
    416 = '4 ч. 04 мин. 51 сек.'
    417 = '1 ч. 42 мин. 02 сек.'
    424 = '2 ч. 10 мин. 08 сек.'
    426 = '3 ч. 01 мин. 57 сек.'
    427 = '3 ч. 07 мин. 30 сек.'
    429 = '4 ч. 07 мин. 01 сек.'
    430 = '2 ч. 00 мин. 35 сек.'
    433 = '2 ч. 07 мин. 24 сек.'
    434 = '2 ч. 26 мин. 02 сек.'
    436 = '2 ч. 01 мин. 30 сек.'
    440 = '2 ч. 01 мин. 21 сек.'
    443 = '1 ч. 59 мин. 03 сек.'
    447 = '5 ч. 43 мин. 08 сек.'
    448 = '2 ч. 08 мин. 54 сек.'
    454 = '7 ч. 28 мин. 05 сек.'
    455 = '9 ч. 22 мин. 05 сек.'
    458 = '3 ч. 08 мин. 28 сек.'
    471 = '3 ч. 02 мин. 11 сек.'
    479 = '2 ч. 07 мин. 39 сек.'
    482 = '7 ч. 17 мин. 04 сек.'
    483 = '13 ч. 01 мин. 41 сек.'
    486 = '2 ч. 23 мин. 05 сек.'
    490 = '5 ч. 03 мин. 21 сек.'
    491 = '8 ч. 43 мин. 02 сек.'
    494 = '8 ч. 03 мин. 47 сек.'
    495 = '10 ч. 03 мин. 38 сек.'
    496 = '5 ч. 06 мин. 06 сек.'
    497 = '2 ч. 06 мин. 51 сек.'
    500 = '2 ч. 05 мин. 37 сек.'
    503 = '26 ч. 02 мин. 17 сек.'
    509 = '11 ч. 15 мин. 05 сек.'
    510 = '2 ч. 02 мин. 40 сек.'
    513 = '3 ч. 01 мин. 15 сек.'
    514 = '2 ч. 46 мин. 05 сек.'
    522 = '2 ч. 30 мин. 03 сек.'
    523 = '4 ч. 48 мин. 04 сек.'
    529 = '4 ч. 10 мин. 08 сек.'
    531 = '5 ч. 04 мин. 13 сек.'
    534 = '8 ч. 08 мин. 30 сек.'
    536 = '3 ч. 53 мин. 01 сек.'
    537 = '4 ч. 00 мин. 29 сек.'
    540 = '6 ч. 29 мин. 03 сек.'
    547 = '4 ч. 09 мин. 22 сек.'
    551 = '3 ч. 25 мин. 04 сек.'
    556 = '5 ч. 46 мин. 06 сек.'
    558 = '3 ч. 10 мин. 02 сек.'
    561 = '4 ч. 09 мин. 28 сек.'
    567 = '4 ч. 04 мин. 25 сек.'
    568 = '6 ч. 09 мин. 01 сек.'
    569 = '3 ч. 47 мин. 01 сек.'
    570 = '9 ч. 03 мин. 15 сек.'
    571 = '5 ч. 08 мин. 07 сек.'
    574 = '4 ч. 17 мин. 06 сек.'
    575 = '10 ч. 03 мин. 25 сек.'
    581 = '5 ч. 02 мин. 14 сек.'
    589 = '5 ч. 51 мин. 05 сек.'
    592 = '3 ч. 55 мин. 09 сек.'
    598 = '13 ч. 48 мин. 09 сек.'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 9).Value = $updates[$row]
}
